$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-9: new forecast-error values in columns B:F,
# and corrected rank values in column G for rows 7-9.
$data = @(
    @(2, 0.4769321400286752, 0.9320513987497278, 1.650466364378836, 1.284704777129297, 1.237927218396544, 14),
    @(3, 0.6221863334236316, 0.9661009109131644, 1.561161447007451, 1.249464464083493, 1.127777791928594, 13),
    @(4, 0.7378115908420845, 0.9692599323336366, 1.482360141223014, 1.217522131717947, 1.011566308962359, 12),
    @(5, 0.7806678668064955, 0.9619392551853039, 2.031973731581155, 1.425473160596563, 1.250913488075273, 11),
    @(6, 0.7978750242939276, 0.8777464132826209, 1.683130036913681, 1.297355015758478, 1.078334869922473, 10),
    @(7, 0.8607107567659571, 0.8607107567659571, 1.054462929432981, 1.026870454065643, 0.5940075024340167, 9),
    @(8, 0.9660439600786086, 1.214125379206839, 2.475674515919399, 1.573427632882872, 1.360485317722336, 6),
    @(9, 1.443592723702684, 1.443592723702684, 2.873932502829098, 1.695267678813319, 1.088558141007014, 3)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}

# Add new row 10 for quarter "Q8" - copy the formatting from the row above
# (bold, centered, bordered label style) and then fill in the values.
$ws.Cells.Item(9, 1).Copy()
$ws.Cells.Item(10, 1).PasteSpecial(-4122)
$ws.Cells.Item(10, 1).Value = "Q8"
$ws.Cells.Item(10, 2).Value = 0.4350014876132097
$ws.Cells.Item(10, 3).Value = 0.4350014876132097
$ws.Cells.Item(10, 4).Value = 0.1892262942257054
$ws.Cells.Item(10, 5).Value = 0.4350014876132097
$ws.Cells.Item(10, 7).Value = 1
